$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset was reprocessed with newly curated dimensions: the
# "municipio-nombre" column's metadata (row 2 concept, row 3 type,
# row 4 value) is updated to reflect that it is now mapped as the
# reference-area dimension instead of a plain measure.
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"
